$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the last three data rows (camasa, blugi, hanorac) from A4:E6 down to J11:N13.
$ws.Range("A4:E6").Cut($ws.Range("J11:N13"))

# Update selection to match the new active range.
$ws.Range("J11:N13").Select()
